# Generate Report for Handback
#
# The "ee06cae0-c551-4ad6-99fb-a3c8ada45c1d.md" file has now been handed
# back (it is in sync with en-US) for both the zh-cn and de-de locales, so
# update the localization-status report accordingly:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" on the Overview sheet as well as
#     on each locale detail sheet.
#   - The "Latest Handback DateTime" for that file is stamped with the
#     handback time for each locale.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn detail sheet -------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("H3").Value = "2016-03-11 16:32:13"

# --- de-de detail sheet -------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("H3").Value = "2016-03-11 16:32:18"

$wb.Save()
